$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITR input data")
$ws.Range("Q4").Value = 14371878
$ws.Range("Q4").NumberFormat = "#,##0"
